$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "MVP" entry that used to live in A1 (the project backlog no
# longer needs a standalone MVP line now that the work is organised into a
# services layer - CurveService will own axes/datasets generation).
$ws.Range("A1").ClearContents()

# Mirror the author's resulting view state: scrolled back to the top of the
# sheet, with the selection left on an empty cell past the data (J21).
$ws.Range("J21").Select()
